$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Header rename: "admin3" -> "admin1" on every sheet (column A header)
# ---------------------------------------------------------------------------
foreach ($wsHeader in $wb.Worksheets) {
    $wsHeader.Cells.Item(1, 1).Value = "admin1"
}

# ---------------------------------------------------------------------------
# Sheet "idp" (sheet1): 11 data rows (2-12) -> 12 data rows (2-13)
# ---------------------------------------------------------------------------
$wsIdp = $wb.Worksheets.Item("idp")
$wsIdp.Rows("13:13").Insert()

$idpData = @(
    ,@(2, "MMR001", "idp", 0.7636430006853658, 0.1576238105360946, 0.07873318877853965, 0)
    ,@(3, "MMR002", "idp", 0.3434720406723505, 0.6160385043040554, 0.01763046118375694, 0.02285899383983707)
    ,@(4, "MMR003", "idp", 0.4831488860607837, 0.4668930084929387, 0.01852596257794513, 0.03143214286833273)
    ,@(5, "MMR004", "idp", 0.2255899434957103, 0.7333886172009896, 0.0283353226488051, 0.01268611665449509)
    ,@(6, "MMR005", "idp", 0.01053189900464389, 0.8701968506970117, 0.02698627938433926, 0.09228497091400525)
    ,@(7, "MMR006", "idp", 0.3666699756698283, 0.5491365975327455, 0.04990483938450553, 0.03428858741292071)
    ,@(8, "MMR007", "idp", 0.665636722421267, 0.2580981627597926, 0.07626511481894054, 0)
    ,@(9, "MMR009", "idp", 0.4379820665480078, 0.4694785473927623, 0.06501086585799726, 0.02752852020123246)
    ,@(10, "MMR011", "idp", 0.4441367798599228, 0.3744393659934109, 0.0644050735709525, 0.1170187805757138)
    ,@(11, "MMR012", "idp", 0.3193362133150582, 0.5546956149147714, 0.1139084598878584, 0.01205971188231214)
    ,@(12, "MMR014", "idp", 0.50562203170395, 0.4512158675159272, 0.0431621007801227, 0)
    ,@(13, "MMR015", "idp", 0.1624601990088694, 0.8132757498888581, 0.02426405110227242, 0)
)

foreach ($row in $idpData) {
    $r = $row[0]
    $wsIdp.Cells.Item($r, 1).Value = $row[1]
    $wsIdp.Cells.Item($r, 2).Value = $row[2]
    $wsIdp.Cells.Item($r, 3).Value = $row[3]
    $wsIdp.Cells.Item($r, 4).Value = $row[4]
    $wsIdp.Cells.Item($r, 5).Value = $row[5]
    $wsIdp.Cells.Item($r, 6).Value = $row[6]
}

# ---------------------------------------------------------------------------
# Sheet "ndsp" (sheet2): 5 data rows (2-6) -> 1 data row (2)
# ---------------------------------------------------------------------------
$wsNdsp = $wb.Worksheets.Item("ndsp")
$wsNdsp.Rows("3:6").Delete()

$ndspData = @(
    ,@(2, "MMR012", "ndsp", 0.3605762564738771, 0.5574700014763981, 0.0763183671935062, 0.00563537485621851)
)

foreach ($row in $ndspData) {
    $r = $row[0]
    $wsNdsp.Cells.Item($r, 1).Value = $row[1]
    $wsNdsp.Cells.Item($r, 2).Value = $row[2]
    $wsNdsp.Cells.Item($r, 3).Value = $row[3]
    $wsNdsp.Cells.Item($r, 4).Value = $row[4]
    $wsNdsp.Cells.Item($r, 5).Value = $row[5]
    $wsNdsp.Cells.Item($r, 6).Value = $row[6]
}

# ---------------------------------------------------------------------------
# Sheet "ocap" (sheet3): 15 data rows (2-16) -> 18 data rows (2-19)
# ---------------------------------------------------------------------------
$wsOcap = $wb.Worksheets.Item("ocap")
$wsOcap.Rows("17:19").Insert()

$ocapData = @(
    ,@(2, "MMR001", "ocap", 0.7236896900952252, 0.2392081879336686, 0.02425551390122181, 0.01284660806988456)
    ,@(3, "MMR002", "ocap", 0.2865024362954217, 0.6836879450220117, 0.02980961868256649, 0)
    ,@(4, "MMR003", "ocap", 0.5184680287850673, 0.4163491503321772, 0.06518282088275558, 0)
    ,@(5, "MMR004", "ocap", 0.3141882758702438, 0.6605474775397754, 0.01981639618841873, 0.005447850401562023)
    ,@(6, "MMR005", "ocap", 0.2568598759208515, 0.6895390258600377, 0.05360109821911088, 0)
    ,@(7, "MMR006", "ocap", 0.6559253601089992, 0.2882136591683537, 0.0558609807226471, 0)
    ,@(8, "MMR007", "ocap", 0.8108905460848412, 0.1891094539151588, 0, 0)
    ,@(9, "MMR008", "ocap", 0.7323658671158727, 0.1779369828077599, 0.0896971500763676, 0)
    ,@(10, "MMR009", "ocap", 0.7714299607478479, 0.2285700392521522, 0, 0)
    ,@(11, "MMR010", "ocap", 0.7596800349328485, 0.1756559775812145, 0.0515602845440673, 0.0131037029418697)
    ,@(12, "MMR011", "ocap", 0.7128105267199996, 0.1403826364618404, 0.08227568805469489, 0.0645311487634652)
    ,@(13, "MMR012", "ocap", 0.4808454996810909, 0.4770717900946745, 0.04208271022423451, 0)
    ,@(14, "MMR013", "ocap", 0.6957578242286829, 0.2236820162678014, 0.08056015950351567, 0)
    ,@(15, "MMR014", "ocap", 0.8475270353110149, 0.1283580034696053, 0.0241149612193799, 0)
    ,@(16, "MMR015", "ocap", 0.4973199688330368, 0.4071193621253184, 0.08437489747912612, 0.01118577156251888)
    ,@(17, "MMR016", "ocap", 0.8660838637861521, 0.1017105875734995, 0.03220554864034835, 0)
    ,@(18, "MMR017", "ocap", 0.9255845165932146, 0.04550899172814852, 0.01723777374349997, 0.01166871793513673)
    ,@(19, "MMR018", "ocap", 0.8751637330175982, 0.06463008325310399, 0.06020618372929779, 0)
)

foreach ($row in $ocapData) {
    $r = $row[0]
    $wsOcap.Cells.Item($r, 1).Value = $row[1]
    $wsOcap.Cells.Item($r, 2).Value = $row[2]
    $wsOcap.Cells.Item($r, 3).Value = $row[3]
    $wsOcap.Cells.Item($r, 4).Value = $row[4]
    $wsOcap.Cells.Item($r, 5).Value = $row[5]
    $wsOcap.Cells.Item($r, 6).Value = $row[6]
}

# ---------------------------------------------------------------------------
# Sheet "ret" (sheet4): 13 data rows (2-14) -> 12 data rows (2-13)
# ---------------------------------------------------------------------------
$wsRet = $wb.Worksheets.Item("ret")
$wsRet.Rows("14:14").Delete()

$retData = @(
    ,@(2, "MMR001", "ret", 0.6109849491331182, 0.3504399005575994, 0.03857515030928242, 0)
    ,@(3, "MMR002", "ret", 0.1389284526796424, 0.7770962716257869, 0.02956559494558583, 0.05440968074898475)
    ,@(4, "MMR003", "ret", 0.5328494875980573, 0.4671505124019426, 0, 0)
    ,@(5, "MMR004", "ret", 0.448486830414843, 0.5272815330456194, 0.02423163653953751, 0)
    ,@(6, "MMR005", "ret", 0.07808646459031106, 0.921913535409689, 0, 0)
    ,@(7, "MMR006", "ret", 0.6683280015228639, 0.2918512553295026, 0.03982074314763352, 0)
    ,@(8, "MMR007", "ret", 0.5370462916961699, 0.4434397190123373, 0.0195139892914926, 0)
    ,@(9, "MMR009", "ret", 0.2290400219628516, 0.5782148690627021, 0.1429641772008985, 0.04978093177354764)
    ,@(10, "MMR011", "ret", 0.5683882818605211, 0.2645865320080541, 0.1565777249880679, 0.01044746114335675)
    ,@(11, "MMR012", "ret", 0.4042210626890028, 0.5226942840398715, 0.06707184162707072, 0.006012811644054996)
    ,@(12, "MMR014", "ret", 0.6195880185893339, 0.1845302652493898, 0.1958817161612764, 0)
    ,@(13, "MMR015", "ret", 0.3339541637595115, 0.5559098893092997, 0.08007043750580958, 0.03006550942537912)
)

foreach ($row in $retData) {
    $r = $row[0]
    $wsRet.Cells.Item($r, 1).Value = $row[1]
    $wsRet.Cells.Item($r, 2).Value = $row[2]
    $wsRet.Cells.Item($r, 3).Value = $row[3]
    $wsRet.Cells.Item($r, 4).Value = $row[4]
    $wsRet.Cells.Item($r, 5).Value = $row[5]
    $wsRet.Cells.Item($r, 6).Value = $row[6]
}

# Re-select A1 on the "idp" sheet to match the original tab-selected view.
$wsIdp.Activate()
$wsIdp.Range("A1").Select()
